$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair below had its Code/Item/Rate/MRP/Qty/Value (columns B-G) rows
# accidentally swapped between two adjacent entries for the same item;
# this restores each row pair by swapping columns B:G back between them.
$rowPairs = @(
    @(183, 184),
    @(279, 280),
    @(313, 314),
    @(317, 318),
    @(346, 347),
    @(355, 356),
    @(372, 373),
    @(379, 380),
    @(382, 383),
    @(389, 390),
    @(400, 401),
    @(421, 422),
    @(431, 432),
    @(457, 458),
    @(536, 537),
    @(581, 582),
    @(586, 587),
    @(590, 591),
    @(593, 594),
    @(599, 600),
    @(604, 605),
    @(687, 688),
    @(709, 710),
    @(715, 716)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("B$r1" + ":G$r1")
    $rng2 = $ws.Range("B$r2" + ":G$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value = $v2
    $rng2.Value = $v1
}

Write-Host "Swapped $($rowPairs.Count) row pairs"
